$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 2
$ws.Range("B2").Value = 4
$ws.Range("B3").Value = 6
$ws.Range("B4").Value = 8
